# Weekly data refresh: a new "Poroto verde" price record (from the most
# recent survey date) is inserted as row 54, pushing all subsequent rows
# (the historical weekly records) down by one. The sheet's used range
# therefore grows from A1:R95 to A1:R96.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 54, shifting row 54..95
# down to 55..96 (Excel copies the formatting from the row above, which
# already carries the date style used by column D).
$ws.Rows.Item(54).EntireRow.Insert()

# Populate the newly inserted row 54 with the new weekly record.
$ws.Range("A54").Value = 7
$ws.Range("B54").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C54").Value = "Ñuble"
$ws.Range("D54").Value = 44810
$ws.Range("E54").Value = 16
$ws.Range("F54").Value = 100112031
$ws.Range("G54").Value = "Poroto verde"
$ws.Range("H54").Value = "Magnum"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 60
$ws.Range("K54").Value = 33000
$ws.Range("L54").Value = 35000
$ws.Range("M54").Value = 34000
$ws.Range("N54").Value = "$/malla 25 kilos"
$ws.Range("O54").Value = "Perú"
$ws.Range("P54").Value = 1360
$ws.Range("Q54").Value = 25
$ws.Range("R54").Value = "Hortaliza"
